$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: new data row appended after existing row 19 (A1:J19 -> A1:J20)
$ws.Cells.Item(20, 1).Value = 45860
$ws.Cells.Item(20, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(20, 2).Value = 6309.6201171875
$ws.Cells.Item(20, 3).Value = 6316.1201171875
$ws.Cells.Item(20, 4).Value = 6281.7099609375
$ws.Cells.Item(20, 5).Value = 6306.60009765625
$ws.Cells.Item(20, 6).Value = 5662040000
$ws.Cells.Item(20, 7).Value = 0.0006375316336257
$ws.Cells.Item(20, 8).Value = 1

$ws.Cells.Item(20, 9).Value = 'promising biotech stock to consider july th best travel stock worth watching july th top small cap stock to follow now july th bailard inc buy share of axon enterprise inc nasdaqaxon airline stock to follow now july th how trump tax law will reshape benefit qa is this the end of the magnificent seven boss stop trying to flatter your worker say leadership expertit doesnt make them like you verkkokauppacom oyj manager transaction nina anttila what separate successful manager from the pack come down to their oneonones here how to do it right today quordle hint and answer for july verkkokauppacom oyj manager transaction pekka litmanen sp signal the u stock rally momentum is waning bloomberg what would you like to leave behind a financial planner guide to family wealth discussion prediction stock that will be worth more than annaly capital year from now top budget gaming phone built for speed top pick under r autohome inc to announce second quarter and interim financial result on july massachusetts financial service co ma sell share of quidelortho corporation nasdaqqdel sbi security co ltd purchase share of fabrinet nysefn apollon wealth management llc acquires share of franklin bsp realty trust inc nysefbrt massachusetts financial service co ma ha million position in royal gold inc nasdaqrgld bitcoin miner iren could surge a sotp valuation peg target hedge fund manager eric jackson say were cooking with gas stock to watch monday verizon bp stellantis block she got the last seat on a flight next to a stranger theyve been married for year zwsoft launch zwcad flex in the u and canada powerful cad performance meet purchasing flexibility cwm llc grows position in lamar advertising company nasdaqlamr top chinese stock to add to your watchlist july th bailard inc sell share of mineral technology inc nysemtx ranger energy service rngr to release earnings on monday hacker exploit microsoft software vulnerability to reportedly target government and businesseswhat to know cwm llc take position in millrose property inc nysemrp acuity inc nyseayi share purchased by cwm llc massachusetts financial service co ma decrease position in select water solution inc nysewttr louisiana state employee retirement system ha stock position in apogee enterprise inc nasdaqapog signaturefd llc lower stock holding in interpublic group of company inc the nyseipg envestnet asset management inc grows holding in yeti holding inc nyseyeti protagonist therapeutic inc nasdaqptgx share sold by massachusetts financial service co ma sbi security co ltd lower position in immunitybio inc nasdaqibrx massachusetts financial service co ma sell share of zurn elkay water solution cor nysezws reynolds consumer product inc nasdaqreyn share sold by massachusetts financial service co ma cwm llc acquires share of spdr sp biotech etf nysearcaxbi edgestream partner lp sell share of the traveler company inc nysetrv sbi security co ltd purchase share of fluor corporation nyseflr st bancorp nasdaqstba rating increased to hold at wall street zen sbi security co ltd ha holding in first trust buywrite income etf nasdaqfthi copa holding sa nysecpa position boosted by massachusetts financial service co ma cemex nysecx cut to hold at wall street zen american asset trust nyseaat upgraded by wall street zen to buy rating massachusetts financial service co ma sell share of kiniksa pharmaceutical international plc nasdaqknsa fg annuity life nysefg upgraded by wall street zen to hold rating brookdale senior living nysebkd downgraded by wall street zen to sell westlake corp nysewlk given consensus rating of moderate buy by brokerage upland software nasdaqupld stock rating upgraded by wall street zen envestnet asset management inc raise stock position in yeti holding inc nyseyeti upland software nasdaqupld rating increased to strongbuy at wall street zen dxc technology nysedxc upgraded by wall street zen to buy rating coupang nysecpng lowered to hold rating by wall street zen dxc technology nysedxc upgraded by wall street zen to buy rating ap business summarybrief at am edt rincon broadcasting buy seven local station from imagicomm rep jonathan l jackson buy international business machine corporation nyseibm share grupo aeroportuario del pacifico nysepac upgraded to buy at wall street zen infineon technology otcmktsifnny stock rating upgraded by wall street zen kennedywilson nysekw upgraded to hold at wall street zen neumora therapeutic inc nasdaqnmra receives consensus rating of hold from analyst insteel industry nyseiiin rating increased to strongbuy at wall street zen world share mostly gain after wall street log a rd straight winning week world share mostly gain after wall street log a rd straight winning week telus digital nysetixt rating lowered to hold at wall street zen nuveen churchill direct lending nysencdl upgraded at wall street zen world share mostly gain after wall street log a rd straight winning week orange otcmktsorany lowered to hold rating by wall street zen wall street zen downgrade irsa inversiones y representaciones nyseirs to buy vista energy nysevist upgraded by wall street zen to hold rating nuveen churchill direct lending nysencdl rating increased to hold at wall street zen telus digital nysetixt downgraded by wall street zen to hold orange otcmktsorany lowered to hold rating by wall street zen gold start the week on the upside a investor weigh trump tariff threat sl green realty nyseslg price target raised to karen kwan sell share of aritzia inc tseatz stock nexxen international nasdaqnexn downgraded by wall street zen to buy nexxen international nasdaqnexn downgraded to buy rating by wall street zen david o watson sell share of apellis pharmaceutical inc nasdaqapls stock victory capital management inc purchase share of sally beauty holding inc nysesbh larson financial group llc acquires share of centerpoint energy inc nysecnp victory capital management inc ha million holding in reddit inc nyserddt larson financial group llc buy share of centerpoint energy inc nysecnp victory capital management inc ha million stock position in highwoods property inc nysehiw microsoft hit with sharepoint attack one version still vulnerable apollon wealth management llc sell share of vaneck oil service etf nysearcaoih apollon wealth management llc grows position in netapp inc nasdaqntap biotechne nasdaqtech and neurogene nasdaqngne financial survey hormel food nysehrl versus china marine food group otcmktscmfo head to head survey sentinelone inc nyse share bought by signaturefd llc sigma planning corp ha million holding in jpmorgan chase co nysejpm cohen capital management inc sell share of jpmorgan chase co nysejpm sigma planning corp ha position in invesco rafi u smallmid etf nasdaqprfz apollon wealth management llc ha position in ppl corporation nyseppl signaturefd llc sell share of illumina inc nasdaqilmn visa inc nysev share sold by edgestream partner lp altman plan dc push to democratize ai economic benefit airline stock to follow now july th blue chip stock worth watching july th doj request for voter data from pa other swing state raise concern among election law expert malcolmjamal warner former cosby show star dy at in costa rica drowning there one thing no one know about trump vladimir tenev sell share of robinhood market inc nasdaqhood stock she got the last seat on a flight next to a stranger theyve been married for year insider selling garrett motion inc nysegtx major shareholder sell share of stock top chinese stock to add to your watchlist july th keith a goldan sell share of syndax pharmaceutical inc nasdaqsndx stock insider selling sintana energy inc cvesei director sell share of stock mineral technology inc nysemtx share sold by bailard inc ford motor nysef upgraded by wall street zen to hold rating performance food group company nysepfgc holding cut by massachusetts financial service co ma echostar corporation nasdaqsats receives consensus recommendation of hold from brokerage the tjx company inc nysetjx receives average pt from brokerage zevra therapeutic nasdaqzvra rating lowered to buy at wall street zen coupang nysecpng downgraded to hold rating by wall street zen irsa inversiones y representaciones nyseirs lowered to buy rating by wall street zen internet initiative japan otcmktsiijiy downgraded to hold rating by wall street zen ofg bancorp nyseofg upgraded by wall street zen to hold rating vista energy nysevist rating increased to hold at wall street zen syrian government start evacuating bedouin family from sweida getting antivaxxers to roll up their sleeve akamai technology nasdaqakam downgraded to hold rating by wall street zen victory capital management inc ha million stock position in yum brand inc nyseyum highwoods property inc nysehiw holding raised by victory capital management inc reviewing hormel food nysehrl and china marine food group otcmktscmfo u envoy double down on support for syria government and criticizes israel intervention larson financial group llc acquires share of ct corporation nysects new parliament seating chart show how labor election win shifted the tectonic plate of federal politics data show more than mbta employee made over k in overtime in cullen frost banker inc invests in sea limited sponsored adr nysese cullen frost banker inc ha holding in howmet aerospace inc nysehwm safehold inc nysesafe given consensus rating of hold by brokerage thomson reuters co nysetri stock position raised by cullen frost banker inc ishares core year usd bond etf nasdaqistb share acquired by cullen frost banker inc apollo global management inc nyseapo share acquired by cullen frost banker inc share in fortis nysefts purchased by cullen frost banker inc costar group inc nasdaqcsgp share acquired by cullen frost banker inc verdence capital advisor llc boost stock position in ishares russell midcap value etf nysearcaiws verdence capital advisor llc ha stake in schwab u midcap etf nysearcaschm verdence capital advisor llc ha stake in medpace holding inc nasdaqmedp verdence capital advisor llc sell share of teck resource ltd nyseteck hillenbrand inc nysehi share acquired by verdence capital advisor llc goldman sachs bdc inc nysegsbd share sold by verdence capital advisor llc verdence capital advisor llc buy new stake in ishares sp value etf nysearcaive tetra tech inc nasdaqttek share acquired by verdence capital advisor llc you can walk between the louvre and the guggenheim in this new art district councilor boston mass and cass drug crisis is worsening the city rat infestation verdence capital advisor llc sell share of target corporation nysetgt bailard inc take position in universal insurance holding inc nyseuve bailard inc ha position in agilent technology inc nysea verdence capital advisor llc buy share of axos financial inc nyseax verdence capital advisor llc purchase new holding in dutch bros inc nysebros verdence capital advisor llc ha holding in millerknoll inc nasdaqmlkn bailard inc grows position in umh property inc nyseumh share in mach natural resource lp nysemnr bought by verdence capital advisor llc victory capital management inc boost stock position in dollar tree inc nasdaqdltr victory capital management inc decrease position in charter communication inc nasdaqchtr larson financial group llc increase position in electronic art inc nasdaqea larson financial group llc raise stock position in airbnb inc nasdaqabnb larson financial group llc sell share of rockwell automation inc nyserok oi glass inc nyseoi receives average recommendation of moderate buy from analyst regret thy name is hawley and murkowski and musk sirius xm holding inc nasdaqsiri share sold by cullen frost banker inc wealthfront adviser llc sell share of ishares russell etf nysearcaiwm wealthfront adviser llc purchase share of the chemours company nysecc wealthfront adviser llc ha holding in first trust senior loan etf nasdaqftsl ast spacemobile inc nasdaqasts share acquired by wealthfront adviser llc pvh corp nysepvh share sold by wealthfront adviser llc emerging market internet etf nysearcaemqq share acquired by wealthfront adviser llc wealthfront adviser llc buy share of f inc nasdaqffiv wealthfront adviser llc boost position in invitation home nyseinvh avantis emerging market equity etf nysearcaavem share purchased by wealthfront adviser llc texas republican want to redraw congressional district in special session inside idaho only maximum security prison where bryan kohberger may spend the rest of his life what to expect a harvard take on the trump administration in federal court vances and marco rubios senate successor avoid gop primary drama independent voter in nyc must count is pakistan moving to join the abraham accord wary democratic establishment keeping new york mamdani at arm length how trump gameplanning for the midterm ai chatbots cocreating a matrix of unreality what i won in a staring contest with a lava heron on san cristobal island vances and marco rubios senate successor avoid gop primary drama bailard inc make new investment in pagaya technology ltd nasdaqpgy share in ge vernova inc nysegev acquired by bailard inc bailard inc sell share of philip morris international inc nysepm promising home improvement stock to research july th bailard inc make new investment in angiodynamics inc nasdaqango best ecommerce stock to watch now july th battery technology stock to research july th nanotechnology stock to consider july th value stock to follow today july th syros pharmaceutical inc nasdaqsyrs receives consensus rating of hold from brokerage massachusetts financial service co ma ha holding in vanguard value etf nysearcavtv arcus bioscience inc nysercus share sold by massachusetts financial service co ma massachusetts financial service co ma sell share of anika therapeutic inc nasdaqanik matson suspends electric vehicle shipment over battery fire concern inside pure fuel partnershipsoverretail strategy sleeper dividend stock to buy for massive upside builder mortgage aid contributing to higher home price morgan stanley say notice on the general meeting of shareholder of the auga group ab entity under restructuring on august to approve the company draft restructuring plan gries financial llc ha million position in procter gamble company the nysepg vanguard growth etf nysearcavug set new year high here why vanguard sp etf nysearcavoo is arrow financial corp largest position gries financial llc ha million stake in bristol myers squibb company nysebmy global bond yield come under pressure obama foundation add outgoing ford foundation president to board of director why opendoor stock is skyrocketing today investopedia jpm high crowding lends risk to stock not involved in ai nucleus security steve carter named a security business innovator hydrograph launch compounding partner program to expand grapheneenhanced thermoplastic klp official deny claim made in civil right lawsuit filed against them lighting precinct reject claim made in civil right lawsuit cornerstone bancorp inc quarterly report june deutsche bank aktiengesellschaft raise american express nyseaxp price target to treasure coast financial planning invests in costco wholesale corporation nasdaqcost share in costco wholesale corporation nasdaqcost purchased by tabor asset management lp costco wholesale corporation nasdaqcost holding trimmed by crestwood advisor group llc im able foundation get donation from firstenergy im able foundation get donation from firstenergy im able foundation get donation from firstenergy rosen law firm encourages vera bradley inc investor to inquire about security class action investigation vra ripple xrp trader who made m earlier this year buy again after month but it not xrp hii summer intern research top use case for ai present finding to company leadership samsung galaxy watch why it the best ai wearable wnba allstar player pay move began well before shirt and sign the pulse business new research suggests rising pr ethic problem with organizational gaslighting best samsung smart tv under r in cnbc daily open investor dismiss trump administration beef with the fed sp hit new high ap technology summarybrief at pm edt delta plane from minneapolis avoids near collision with b bomber iran sends a rocket designed to carry satellite into a suborbital test flight iran sends a rocket designed to carry satellite into a suborbital test flight wex tap ai for faster fsa reimbursement ap business summarybrief at pm edt biogen bet big on billion expansion to boost u drug manufacturing retail speculation is back with a vengeance heat pump industry on the rise global market to expand significantly by senior connection receives innovative achievement award for care express weekly analyst rating update for metlife met hazlett burt watson inc boost stock position in schwab u dividend equity etf nysearcaschd hazlett burt watson inc grows holding in salesforce inc nysecrm verus financial partner inc raise stock position in schwab u dividend equity etf nysearcaschd western montana mental health center wmmhc data breach expose personal information murphy law firm investigates legal claim appeal court order new trial for man convicted in etan patz case animal science expertise honored at international conference it end with u insurer file lawsuit say it doesnt have to pay justin baldonis legal fee jonas brother sell portion of music catalog to dad company i dont know that anybody can love your music more silicon lab announces second quarter earnings webcast head to head review ge aerospace nysege v northrop grumman nysenoc van strum towne inc boost holding in mcdonalds corporation nysemcd analyzing storage computer otcmktssoso teradata nysetdc financial review crh crh v it rival gauzy nasdaqgauz versus kla nasdaqklac head to head comparison netflix inc nasdaqnflx share sold by lewis asset management llc replimune group nasdaqrepl see large volume increase should you buy pacific wealth management sell share of abbott laboratory nyseabt walker dunlop nysewd v rocket company nyserkt headtohead survey zegona communication lonzeg hit new year high following analyst upgrade dwave quantum stock is trading higher monday whats going on rogers returning to bargaining table tomorrow usw local member call for real action toward a fair deal southwest airline announces start date for assigned seating cnn democrat try to dissect themselves with autopsy show theyre riddled with disease vista capital partner inc trim stock position in nike inc nysenke van strum towne inc ha million position in pepsico inc nasdaqpep van strum towne inc sell share of ishares core sp etf nysearcaivv integrity investment advisor llc acquires share of ishares core sp etf nysearcaivv unit apartment project proposed for virginia ave in hollywood integrity investment advisor llc buy new stake in ishares core sp etf nysearcaivv astronomer ceo kiss cam controversy sparked over million in prediction market bet on his ouster whataburger saddle up with the new bacon wrangler double a burger a big and bold a the lone star state rosen a ranked and leading law firm encourages petco health and wellness company inc investor to secure counsel before important deadline in security class action woof critical comparison finwise bancorp nasdaqfinw and texas capital bancshares nasdaqtcbi analyzing kuehne nagel international khngy it competitor rd lewis holding inc acquires share of walmart inc nysewmt comparing rayonier advanced material nyseryam and suzano nysesuz comparing osisko mining otcmktsobnnf aris mining nysearmn analyzing elme community nyseelme sun community nysesui marinemax hzo v it competitor headtohead analysis contrasting williamssonoma wsm and it competitor oneascent wealth management llc increase stake in walmart inc nysewmt reviewing ivanhoe mine ivpaf and it competitor contrasting senstar technology nasdaqsnt applied dna science nasdaqapdn ivanhoe mine ivpaf and it competitor financial comparison headtohead analysis suzano nysesuz versus rayonier advanced material nyseryam lewis asset management llc increase stake in cisco system inc nasdaqcsco unifirst nyseunf and vestis nysevsts financial analysis virgin galactic nysespce v safran otcmktssafry critical contrast amazon ups quantum game with new qubit qpu eye commercial future perception that tariff will be pulled back if market react is alive and well say cfrs patterson skip brittenham hollywood lawyer to the star dy at rogers announces pricing of cash tender offer for canadian dollar debt security rep gilbert ray cisneros jr buy northrop grumman corporation nysenoc stock rep gilbert ray cisneros jr sell prologis inc nysepld stock rep robert bresnahan jr sell off share of arista network inc nyseanet rep gilbert ray cisneros jr purchase share of okta inc nasdaqokta aire advisor llc cut stock position in abbvie inc nyseabbv lauterbach financial advisor llc invests in vanguard sp etf nysearcavoo minnesota lawmaker to resign after being convicted of felony burglary research offer link between burn pit smoke and serious brain injury gop megabills final score t in red ink and million kicked off health insurance cbo say gop megabills final score t in red ink and million kicked off health insurance cbo say trump wade into bryan kohberger murder case demand answer on why he did it there are no explanation sentencing hearing set for exkentucky officer convicted in breonna taylor raid treasure coast financial planning purchase share of costco wholesale corporation nasdaqcost ge vernova gev versus it peer critical contrast tabor asset management lp make new investment in costco wholesale corporation nasdaqcost wealth dimension group ltd ha stake in cocacola company the nyseko trump team throw nyc mayor eric adam under the bus the public deserves to know dcfs withholds detail in illinois girl death miller financial service llc sell share of vanguard value etf nysearcavtv maine need to brace for natural disaster opinion maine need to brace for natural disaster opinion maine need to brace for natural disaster opinion appeal court order new trial for man convicted in etan patz case hunter biden go off on george clooney f him and everybody around him activeduty marine will be leaving los angeles pentagon confirms encountering a bear on the trail yieldmax ultra option income strategy etf nysearcaulty see strong trading volume still a buy contrasting cosan csan it rival rollins nyserol versus limbach nasdaqlmb critical comparison southwest airline announces start date for assigned seating rogers returning to bargaining table tomorrow usw local member call for real action toward a fair deal bruker nasdaqbrkr share gap down whats next irobot nasdaqirbt v robogroup tek otcmktsrobof head to head analysis headtohead contrast kuehne nagel international khngy and it peer srh advisor llc cut holding in the charles schwab corporation nyseschw lewis asset management llc sell share of pfizer inc nysepfe limoneira nasdaqlmnr and forafric global nasdaqafri headtohead review analyzing samsonite group smsey it competitor senstar technology nasdaqsnt and applied dna science nasdaqapdn head to head survey malcolmjamal warner cosby show star dy of drowning at okta nasdaqokta stock unloaded rep robert bresnahan jr travis county commissioner margaret gomez to retire after more than year of service resecurity and braly insurance group announce strategic partnership to drive innovation in cybersecurity and insurance risk management pb and npr are generally unbiased independent of government propaganda and provide key benefit to u democracy this analyst just returned from ukraine frontlines here what he saw judge weighs reality of trump ideological deportation policy a activist crackdown trial end judge weighs reality of trump ideological deportation policy a activist crackdown trial end judge weighs reality of trump ideological deportation policy a activist crackdown trial end top metaverse stock to follow today july st judge weighs reality of trump ideological deportation policy a activist crackdown trial end judge weighs reality of trump ideological deportation policy a activist crackdown trial end lowering the voting age is a call to trust young people with democracy likely new owner of paramount pledge to snuff out bias at cbs lawyer tell fcc that cbs news will practice unbiased journalism judge weighs reality of trump ideological deportation policy a activist crackdown trial end no fk you hunter julio rosa wallop hunter biden for whining about illegals not cleaning his toilet the day niagara fall ran dry here what it looked like highpoint advisor group llc boost stock holding in johnson johnson nysejnj highpoint advisor group llc purchase share of uber technology inc nyseuber hotel stock to add to your watchlist july st highpoint advisor group llc ha million stake in bank of america corporation nysebac promising g stock to add to your watchlist july st peterson wealth service purchase share of johnson johnson nysejnj promising digital medium stock to follow today july st best gold stock to watch today july st verizon communication nysevz release quarterly earnings result beat expectation by eps exxon mobil nysexom braskem nysebak headtohead contrast financial stock to keep an eye on july st johnson johnson nysejnj stock holding decreased by peregrine asset adviser inc this pb station in boston sign is going viral after trump defunded them mark herr communication issue statement regarding the sentencing of anne pramaggiore mastercard incorporated nysema share sold by gries financial llc trump say make indian great again year ago he tried to buy and move the cleveland indian hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance visa inc nysev stock position lowered by miller financial service llc arrow financial corp reduces stock holding in the home depot inc nysehd peterson wealth service purchase share of visa inc nysev jw cole advisor inc purchase share of oracle corporation nyseorcl harvard argues against trump admins funding freeze in court hearing future generation global limited asxfgg insider geoff wilson purchase share d boral capital reiterates buy rating for biomea fusion nasdaqbmea pepsi is betting a popular health trend can get more people to drink it soda america golden age is over letter epstein ghost is haunting trump and america enemy could summon more how did being a sport fan get so expensive wealth dimension group ltd boost position in well fargo company nysewfc environmentalist lawsuit to halt alligator alcatraz filed in wrong court florida official say trump official who hold more than one job malcolmjamal warner actor who starred a theo in the cosby show dead at young kim draw another democratic challenger in california hazlett burt watson inc acquires share of caterpillar inc nysecat hazlett burt watson inc ha million stock holding in lockheed martin corporation nyselmt gop house majority shrink after big beautiful victory dhs release video showing attempted robbery of offduty federal officer in nyc update dare bioscience to host august webinar the dare to playtm difference the winch advisory service llc ha million holding in salesforce inc nysecrm clinton email investigation omitted crucial piece of evidence declassified report find sei select international equity etf nasdaqseie see strong trading volume should you buy head to head comparison hydrofarm holding group nasdaqhyfm one bio otcmktsonbi concord medical service nyseccm versus brookdale senior living nysebkd critical contrast'
$ws.Cells.Item(20, 10).Value = 'promising biotech stock consider best travel stock worth watching small cap stock follow bailard buy axon enterprise nasdaqaxon airline stock follow trump tax reshape benefit qa end magnificent seven bos stop flatter worker leadership expertit doesnt verkkokauppacom oyj manager transaction nina anttila separate successful manager pack come oneonones quordle hint answer verkkokauppacom oyj manager transaction pekka litmanen sp signal stock rally momentum waning bloomberg behind financial planner guide family wealth discussion prediction stock worth annaly budget gaming phone built speed pick r autohome announce quarter interim financial massachusetts financial service co sell quidelortho corporation nasdaqqdel sbi security co ltd purchase fabrinet nysefn apollon wealth acquires franklin bsp realty trust nysefbrt massachusetts financial service co million royal gold nasdaqrgld bitcoin miner iren surge sotp valuation peg target hedge fund manager eric jackson cooking gas stock watch monday verizon bp stellantis block got seat flight next stranger theyve married zwsoft launch zwcad flex canada powerful cad performance meet purchasing flexibility cwm grows lamar advertising company nasdaqlamr chinese stock add watchlist bailard sell mineral technology nysemtx ranger energy service rngr release earnings monday hacker exploit microsoft software vulnerability reportedly target government businesseswhat cwm take millrose property nysemrp acuity nyseayi purchased cwm massachusetts financial service co decrease select water solution nysewttr louisiana employee retirement system stock apogee enterprise nasdaqapog signaturefd lower stock interpublic company nyseipg envestnet asset grows yeti nyseyeti protagonist therapeutic nasdaqptgx massachusetts financial service co sbi security co ltd lower immunitybio nasdaqibrx massachusetts financial service co sell zurn elkay water solution cor nysezws reynolds consumer product nasdaqreyn massachusetts financial service co cwm acquires spdr sp biotech nysearcaxbi edgestream partner lp sell traveler company nysetrv sbi security co ltd purchase fluor corporation nyseflr st bancorp nasdaqstba rating increased wall street zen sbi security co ltd trust buywrite income nasdaqfthi copa sa nysecpa boosted massachusetts financial service co cemex nysecx cut wall street zen american asset trust nyseaat upgraded wall street zen buy rating massachusetts financial service co sell kiniksa pharmaceutical international plc nasdaqknsa fg annuity life nysefg upgraded wall street zen rating brookdale senior living nysebkd downgraded wall street zen sell westlake corp nysewlk consensus rating moderate buy brokerage upland software nasdaqupld stock rating upgraded wall street zen envestnet asset raise stock yeti nyseyeti upland software nasdaqupld rating increased strongbuy wall street zen dxc technology nysedxc upgraded wall street zen buy rating coupang nysecpng lowered rating wall street zen dxc technology nysedxc upgraded wall street zen buy rating business rincon broadcasting buy seven local station imagicomm rep jonathan l jackson buy international business machine corporation nyseibm grupo aeroportuario del pacifico nysepac upgraded buy wall street zen infineon technology otcmktsifnny stock rating upgraded wall street zen kennedywilson nysekw upgraded wall street zen neumora therapeutic nasdaqnmra receives consensus rating analyst insteel industry nyseiiin rating increased strongbuy wall street zen world mostly gain wall street log rd straight winning world mostly gain wall street log rd straight winning telus digital nysetixt rating lowered wall street zen nuveen churchill direct lending nysencdl upgraded wall street zen world mostly gain wall street log rd straight winning orange otcmktsorany lowered rating wall street zen wall street zen downgrade irsa inversiones representaciones nyseirs buy vista energy nysevist upgraded wall street zen rating nuveen churchill direct lending nysencdl rating increased wall street zen telus digital nysetixt downgraded wall street zen orange otcmktsorany lowered rating wall street zen gold start upside investor weigh trump tariff threat sl green realty nyseslg price target raised karen kwan sell aritzia tseatz stock nexxen international nasdaqnexn downgraded wall street zen buy nexxen international nasdaqnexn downgraded buy rating wall street zen david watson sell apellis pharmaceutical nasdaqapls stock victory purchase sally beauty nysesbh larson financial acquires centerpoint energy nysecnp victory million reddit nyserddt larson financial buy centerpoint energy nysecnp victory million stock highwoods property nysehiw microsoft hit sharepoint attack version vulnerable apollon wealth sell vaneck oil service nysearcaoih apollon wealth grows netapp nasdaqntap biotechne nasdaqtech neurogene nasdaqngne financial survey hormel nysehrl versus china marine otcmktscmfo head head survey sentinelone nyse bought signaturefd sigma planning corp million jpmorgan chase co nysejpm cohen sell jpmorgan chase co nysejpm sigma planning corp invesco rafi smallmid nasdaqprfz apollon wealth ppl corporation nyseppl signaturefd sell illumina nasdaqilmn visa nysev edgestream partner lp altman plan dc push democratize ai economic benefit airline stock follow blue chip stock worth watching doj request voter data pa swing raise concern among election expert malcolmjamal warner cosby star dy costa rica drowning trump vladimir tenev sell robinhood market nasdaqhood stock got seat flight next stranger theyve married insider selling garrett motion nysegtx major shareholder sell stock chinese stock add watchlist keith goldan sell syndax pharmaceutical nasdaqsndx stock insider selling sintana energy cvesei director sell stock mineral technology nysemtx bailard ford nysef upgraded wall street zen rating performance company nysepfgc cut massachusetts financial service co echostar corporation nasdaqsats receives consensus recommendation brokerage tjx company nysetjx receives average pt brokerage zevra therapeutic nasdaqzvra rating lowered buy wall street zen coupang nysecpng downgraded rating wall street zen irsa inversiones representaciones nyseirs lowered buy rating wall street zen internet initiative japan otcmktsiijiy downgraded rating wall street zen ofg bancorp nyseofg upgraded wall street zen rating vista energy nysevist rating increased wall street zen syrian government start evacuating bedouin family sweida getting antivaxxers roll sleeve akamai technology nasdaqakam downgraded rating wall street zen victory million stock yum brand nyseyum highwoods property nysehiw raised victory reviewing hormel nysehrl china marine otcmktscmfo envoy double support syria government criticizes israel intervention larson financial acquires ct corporation nysects parliament seating chart labor election shifted tectonic plate federal politics data mbta employee k overtime cullen frost banker invests sea limited sponsored adr nysese cullen frost banker howmet aerospace nysehwm safehold nysesafe consensus rating brokerage thomson reuters co nysetri stock raised cullen frost banker ishares core usd bond nasdaqistb acquired cullen frost banker apollo global nyseapo acquired cullen frost banker fortis nysefts purchased cullen frost banker costar nasdaqcsgp acquired cullen frost banker verdence advisor boost stock ishares russell midcap value nysearcaiws verdence advisor stake schwab midcap nysearcaschm verdence advisor stake medpace nasdaqmedp verdence advisor sell teck resource ltd nyseteck hillenbrand nysehi acquired verdence advisor goldman sachs bdc nysegsbd verdence advisor verdence advisor buy stake ishares sp value nysearcaive tetra tech nasdaqttek acquired verdence advisor walk louvre guggenheim art district councilor boston mass cass drug crisis worsening rat infestation verdence advisor sell target corporation nysetgt bailard take universal insurance nyseuve bailard agilent technology nysea verdence advisor buy axos financial nyseax verdence advisor purchase dutch bros nysebros verdence advisor millerknoll nasdaqmlkn bailard grows umh property nyseumh mach natural resource lp nysemnr bought verdence advisor victory boost stock dollar tree nasdaqdltr victory decrease charter communication nasdaqchtr larson financial increase electronic art nasdaqea larson financial raise stock airbnb nasdaqabnb larson financial sell rockwell automation nyserok oi glass nyseoi receives average recommendation moderate buy analyst regret thy name hawley murkowski musk sirius xm nasdaqsiri cullen frost banker wealthfront adviser sell ishares russell nysearcaiwm wealthfront adviser purchase chemours company nysecc wealthfront adviser trust senior loan nasdaqftsl ast spacemobile nasdaqasts acquired wealthfront adviser pvh corp nysepvh wealthfront adviser emerging market internet nysearcaemqq acquired wealthfront adviser wealthfront adviser buy f nasdaqffiv wealthfront adviser boost invitation nyseinvh avantis emerging market equity nysearcaavem purchased wealthfront adviser texas republican redraw congressional district special session inside idaho maximum security prison bryan kohberger spend rest life expect harvard take trump administration federal vances marco rubios senate successor avoid gop primary drama independent voter nyc must count pakistan moving join abraham wary democratic establishment keeping york mamdani arm length trump gameplanning midterm ai chatbots cocreating matrix unreality staring contest lava heron san cristobal island vances marco rubios senate successor avoid gop primary drama bailard investment pagaya technology ltd nasdaqpgy ge vernova nysegev acquired bailard bailard sell philip morris international nysepm promising improvement stock research bailard investment angiodynamics nasdaqango best ecommerce stock watch battery technology stock research nanotechnology stock consider value stock follow syros pharmaceutical nasdaqsyrs receives consensus rating brokerage massachusetts financial service co vanguard value nysearcavtv arcus bioscience nysercus massachusetts financial service co massachusetts financial service co sell anika therapeutic nasdaqanik matson suspends electric vehicle shipment battery fire concern inside pure fuel partnershipsoverretail strategy sleeper dividend stock buy massive upside builder mortgage aid contributing higher price morgan stanley notice general meeting shareholder auga ab entity restructuring august approve company draft restructuring plan gries financial million procter gamble company nysepg vanguard growth nysearcavug set high vanguard sp nysearcavoo arrow financial corp largest gries financial million stake bristol myers squibb company nysebmy global bond yield come pressure obama foundation add outgoing ford foundation president board director opendoor stock skyrocketing investopedia jpm high crowding lends risk stock involved ai nucleus security steve carter named security business innovator hydrograph launch compounding partner program expand grapheneenhanced thermoplastic klp deny claim civil lawsuit filed lighting precinct reject claim civil lawsuit cornerstone bancorp quarterly report deutsche bank aktiengesellschaft raise american express nyseaxp price target treasure coast financial planning invests costco wholesale corporation nasdaqcost costco wholesale corporation nasdaqcost purchased tabor asset lp costco wholesale corporation nasdaqcost trimmed crestwood advisor im able foundation donation firstenergy im able foundation donation firstenergy im able foundation donation firstenergy rosen firm encourages vera bradley investor inquire security class action investigation vra ripple xrp trader buy month xrp hii summer intern research ai present finding company leadership samsung galaxy watch best ai wearable wnba allstar player pay move began shirt sign pulse business research suggests rising pr ethic problem organizational gaslighting best samsung smart tv r cnbc open investor dismiss trump administration beef fed sp hit high technology pm delta plane minneapolis avoids near collision b bomber iran sends rocket designed carry satellite suborbital test flight iran sends rocket designed carry satellite suborbital test flight wex tap ai faster fsa reimbursement business pm biogen bet billion expansion boost drug manufacturing retail speculation vengeance heat pump industry rise global market expand significantly senior connection receives innovative achievement award care express weekly analyst rating update metlife met hazlett burt watson boost stock schwab dividend equity nysearcaschd hazlett burt watson grows salesforce nysecrm verus financial partner raise stock schwab dividend equity nysearcaschd western montana mental health center wmmhc data breach expose personal information murphy firm investigates legal claim appeal order trial man convicted etan patz animal science expertise honored international conference end insurer file lawsuit doesnt pay justin baldonis legal fee jonas brother sell portion music catalog dad company dont anybody love music silicon lab announces quarter earnings webcast head head review ge aerospace nysege v northrop grumman nysenoc van strum towne boost mcdonalds corporation nysemcd analyzing storage computer otcmktssoso teradata nysetdc financial review crh crh v rival gauzy nasdaqgauz versus kla nasdaqklac head head comparison netflix nasdaqnflx lewis asset replimune nasdaqrepl large volume increase buy pacific wealth sell abbott laboratory nyseabt walker dunlop nysewd v rocket company nyserkt headtohead survey zegona communication lonzeg hit high following analyst upgrade dwave quantum stock trading higher monday whats rogers returning bargaining table tomorrow usw local call real action fair deal southwest airline announces start date assigned seating cnn democrat try dissect autopsy theyre riddled disease vista partner trim stock nike nysenke van strum towne million pepsico nasdaqpep van strum towne sell ishares core sp nysearcaivv integrity investment advisor acquires ishares core sp nysearcaivv unit apartment project proposed virginia ave hollywood integrity investment advisor buy stake ishares core sp nysearcaivv astronomer ceo kiss cam controversy sparked million prediction market bet ouster whataburger saddle bacon wrangler double burger bold lone star rosen ranked leading firm encourages petco health wellness company investor secure counsel important deadline security class action woof critical comparison finwise bancorp nasdaqfinw texas bancshares nasdaqtcbi analyzing kuehne nagel international khngy competitor rd lewis acquires walmart nysewmt comparing rayonier advanced material nyseryam suzano nysesuz comparing osisko mining otcmktsobnnf aris mining nysearmn analyzing elme community nyseelme sun community nysesui marinemax hzo v competitor headtohead analysis contrasting williamssonoma wsm competitor oneascent wealth increase stake walmart nysewmt reviewing ivanhoe mine ivpaf competitor contrasting senstar technology nasdaqsnt applied dna science nasdaqapdn ivanhoe mine ivpaf competitor financial comparison headtohead analysis suzano nysesuz versus rayonier advanced material nyseryam lewis asset increase stake cisco system nasdaqcsco unifirst nyseunf vestis nysevsts financial analysis virgin galactic nysespce v safran otcmktssafry critical contrast amazon ups quantum game qubit qpu eye commercial future perception tariff pulled market react alive cfrs patterson skip brittenham hollywood lawyer star dy rogers announces pricing cash tender offer canadian dollar debt security rep gilbert ray cisneros jr buy northrop grumman corporation nysenoc stock rep gilbert ray cisneros jr sell prologis nysepld stock rep robert bresnahan jr sell arista network nyseanet rep gilbert ray cisneros jr purchase okta nasdaqokta aire advisor cut stock abbvie nyseabbv lauterbach financial advisor invests vanguard sp nysearcavoo minnesota lawmaker resign convicted felony burglary research offer link burn pit smoke serious brain injury gop megabills final score red ink million kicked health insurance cbo gop megabills final score red ink million kicked health insurance cbo trump wade bryan kohberger murder demand answer explanation sentencing hearing set exkentucky officer convicted breonna taylor raid treasure coast financial planning purchase costco wholesale corporation nasdaqcost ge vernova gev versus peer critical contrast tabor asset lp investment costco wholesale corporation nasdaqcost wealth dimension ltd stake cocacola company nyseko trump team throw nyc mayor eric adam bus public deserves dcfs withholds detail illinois girl death miller financial service sell vanguard value nysearcavtv maine brace natural disaster opinion maine brace natural disaster opinion maine brace natural disaster opinion appeal order trial man convicted etan patz hunter biden george clooney f everybody around activeduty marine leaving los angeles pentagon confirms encountering bear trail yieldmax ultra option income strategy nysearcaulty strong trading volume buy contrasting cosan csan rival rollins nyserol versus limbach nasdaqlmb critical comparison southwest airline announces start date assigned seating rogers returning bargaining table tomorrow usw local call real action fair deal bruker nasdaqbrkr gap whats next irobot nasdaqirbt v robogroup tek otcmktsrobof head head analysis headtohead contrast kuehne nagel international khngy peer srh advisor cut charles schwab corporation nyseschw lewis asset sell pfizer nysepfe limoneira nasdaqlmnr forafric global nasdaqafri headtohead review analyzing samsonite smsey competitor senstar technology nasdaqsnt applied dna science nasdaqapdn head head survey malcolmjamal warner cosby star dy drowning okta nasdaqokta stock unloaded rep robert bresnahan jr travis county commissioner margaret gomez retire service resecurity braly insurance announce strategic partnership drive innovation cybersecurity insurance risk pb npr generally unbiased independent government propaganda provide benefit democracy analyst returned ukraine frontlines saw judge weighs reality trump ideological deportation policy activist crackdown trial end judge weighs reality trump ideological deportation policy activist crackdown trial end judge weighs reality trump ideological deportation policy activist crackdown trial end metaverse stock follow st judge weighs reality trump ideological deportation policy activist crackdown trial end judge weighs reality trump ideological deportation policy activist crackdown trial end lowering voting age call trust young people democracy likely owner paramount pledge snuff bias cbs lawyer fcc cbs practice unbiased journalism judge weighs reality trump ideological deportation policy activist crackdown trial end fk hunter julio rosa wallop hunter biden whining illegals cleaning toilet niagara fall ran dry looked highpoint advisor boost stock johnson johnson nysejnj highpoint advisor purchase uber technology nyseuber hotel stock add watchlist st highpoint advisor million stake bank america corporation nysebac promising g stock add watchlist st peterson wealth service purchase johnson johnson nysejnj promising digital medium stock follow st best gold stock watch st verizon communication nysevz release quarterly earnings beat expectation eps exxon mobil nysexom braskem nysebak headtohead contrast financial stock eye st johnson johnson nysejnj stock decreased peregrine asset adviser pb station boston sign viral trump defunded mark herr communication issue statement regarding sentencing anne pramaggiore mastercard incorporated nysema gries financial trump indian great ago tried buy move cleveland indian hunter biden suggests ambien contributed joe bidens poor debate performance hunter biden suggests ambien contributed joe bidens poor debate performance hunter biden suggests ambien contributed joe bidens poor debate performance hunter biden suggests ambien contributed joe bidens poor debate performance hunter biden suggests ambien contributed joe bidens poor debate performance hunter biden suggests ambien contributed joe bidens poor debate performance visa nysev stock lowered miller financial service arrow financial corp reduces stock depot nysehd peterson wealth service purchase visa nysev jw cole advisor purchase oracle corporation nyseorcl harvard argues trump admins funding freeze hearing future generation global limited asxfgg insider geoff wilson purchase boral reiterates buy rating biomea fusion nasdaqbmea pepsi betting popular health trend people drink soda america golden age epstein ghost haunting trump america enemy summon sport fan expensive wealth dimension ltd boost fargo company nysewfc environmentalist lawsuit halt alligator alcatraz filed florida trump job malcolmjamal warner actor starred theo cosby dead young kim draw another democratic challenger california hazlett burt watson acquires caterpillar nysecat hazlett burt watson million stock lockheed martin corporation nyselmt gop house majority shrink victory dhs release video showing attempted robbery offduty federal officer nyc update dare bioscience host august webinar dare playtm difference winch advisory service million salesforce nysecrm clinton email investigation omitted crucial piece evidence declassified report find sei select international equity nasdaqseie strong trading volume buy head head comparison hydrofarm nasdaqhyfm bio otcmktsonbi concord medical service nyseccm versus brookdale senior living nysebkd critical contrast'
